$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename BOM item labels in column A to reflect updated fastener naming
$ws.Range("A13").Value = "M6x15 Bolt"
$ws.Range("A14").Value = "M6x20 Bolt"
$ws.Range("A19").Value = "M4x14 Bolt"
$ws.Range("A21").Value = "M3x25 Bolt"
$ws.Range("A22").Value = "M3x10 Pan Head Bolt"
$ws.Range("A17").Value = "M5x40 Bolt"

# Update the active sheet view: show formulas and move the selection
$ws.Activate()
$excel.ActiveWindow.DisplayFormulas = $true
$ws.Range("C18").Select()
